$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Version: 0.1.1 -> 0.2.0
$ws.Range("B3").Value = "0.2.0"

# Date: 2023-10-20T07:19:33+00:00 -> 2023-10-20T08:59:58+00:00
$ws.Range("B8").Value = "2023-10-20T08:59:58+00:00"

# Insert a new "Jurisdiction" row right after "Contact" (row 10) and before
# "Description" (old row 11), pushing Description/Purpose/Copyright/Immutable
# down by one row.
$ws.Rows.Item(11).Insert()
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = "iso:code:3166:FR"

# Make the new row match the look (borders/fill/wrap) of the other data rows.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

Write-Output "done"
